# Update HotStock_Top20 rankings (A2:C21) with new stock names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "实达集团"
$ws.Range("B2").Value = "航天发展"
$ws.Range("C2").Value = "实达集团"
$ws.Range("B3").Value = "蓝色光标"
$ws.Range("A4").Value = "蓝色光标"
$ws.Range("C4").Value = "榕基软件"
$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "实达集团"
$ws.Range("A6").Value = "光库科技"
$ws.Range("B6").Value = "榕基软件"
$ws.Range("C6").Value = "摩尔线程"
$ws.Range("A7").Value = "榕基软件"
$ws.Range("B7").Value = "雷科防务"
$ws.Range("C7").Value = "雷科防务"
$ws.Range("A8").Value = "特发信息"
$ws.Range("B8").Value = "省广集团"
$ws.Range("C8").Value = "天际股份"
$ws.Range("A9").Value = "天际股份"
$ws.Range("B9").Value = "长城军工"
$ws.Range("C9").Value = "特发信息"
$ws.Range("A10").Value = "雷科防务"
$ws.Range("B10").Value = "三六零"
$ws.Range("C10").Value = "蓝色光标"
$ws.Range("A11").Value = "省广集团"
$ws.Range("B11").Value = "天际股份"
$ws.Range("C11").Value = "合富中国"
$ws.Range("A12").Value = "合富中国"
$ws.Range("B12").Value = "光库科技"
$ws.Range("C12").Value = "国风新材"
$ws.Range("A13").Value = "航天动力"
$ws.Range("B13").Value = "久其软件"
$ws.Range("C13").Value = "海安集团"
$ws.Range("A14").Value = "三六零"
$ws.Range("B14").Value = "合富中国"
$ws.Range("C14").Value = "海南海药"
$ws.Range("A15").Value = "新 华 都"
$ws.Range("B15").Value = "海南海药"
$ws.Range("C15").Value = "光库科技"
$ws.Range("A16").Value = "海南海药"
$ws.Range("B16").Value = "特发信息"
$ws.Range("C16").Value = "省广集团"
$ws.Range("A17").Value = "久其软件"
$ws.Range("B17").Value = "永鼎股份"
$ws.Range("C17").Value = "国晟科技"
$ws.Range("A18").Value = "国风新材"
$ws.Range("B18").Value = "多氟多"
$ws.Range("C18").Value = "中水渔业"
$ws.Range("A19").Value = "易点天下"
$ws.Range("B19").Value = "国风新材"
$ws.Range("C19").Value = "航天动力"
$ws.Range("A20").Value = "中水渔业"
$ws.Range("B20").Value = "中水渔业"
$ws.Range("C20").Value = "久其软件"
$ws.Range("A21").Value = "长城军工"
$ws.Range("B21").Value = "航天动力"
$ws.Range("C21").Value = "欢瑞世纪"
